# Rename category labels in column E from "Language of Other ..." to
# "Languages of Other ..." (pluralize "Language" -> "Languages").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)  # Column E
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "Language of Other*") {
        $cell.Value = $val -replace "^Language of Other", "Languages of Other"
    }
}
